# The data rows 2-10 (columns D, M, Q, S, T) have been permuted among each
# other (the rest of the row content - A,B,C,E-L,N,O,P,R - is identical for
# every row so it is left untouched). Capture the current values first, then
# write them back out in their new row positions so no values are lost while
# being overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values for the columns that move.
$rows = 2..10
$snapshot = @{}
foreach ($r in $rows) {
    $snapshot[$r] = @{
        D = $ws.Cells.Item($r, 4).Value2
        M = $ws.Cells.Item($r, 13).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

# Mapping of source row (old) -> destination row (new) for the row contents.
$mapping = @{
    2 = 4
    3 = 10
    4 = 6
    5 = 8
    6 = 2
    7 = 5
    8 = 9
    9 = 3
    10 = 7
}

foreach ($oldRow in $mapping.Keys) {
    $newRow = $mapping[$oldRow]
    $data = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 13).Value = $data.M
    $ws.Cells.Item($newRow, 17).Value = $data.Q
    $ws.Cells.Item($newRow, 19).Value = $data.S
    $ws.Cells.Item($newRow, 20).Value = $data.T
}
